$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order rows appended to the bottom of the sheet (rows 14-19).
# Numeric-looking values are written as text (matching the source
# workbook's convention of storing quantities/prices as strings),
# using a leading apostrophe so Excel keeps them as text instead of
# auto-converting to numbers. Empty SKU cells are written the same
# way so the cell exists (empty string) rather than being left blank,
# matching the existing empty SKU cells above them.

# Row 14 - Pineapple - Fresh
$ws.Range("A14").Value = "'"
$ws.Range("B14").Value = "Pineapple - Fresh"
$ws.Range("C14").Value = "'3"
$ws.Range("D14").Value = "'22.50"
$ws.Range("E14").Value = "'67.50"

# Row 15 - Mushroom - White Sliced
$ws.Range("A15").Value = "'"
$ws.Range("B15").Value = "Mushroom - White Sliced"
$ws.Range("C15").Value = "'2"
$ws.Range("D15").Value = "'22.50"
$ws.Range("E15").Value = "'45.00"

# Row 16 - Scallion - Fresh (SKU: Green Onion)
$ws.Range("A16").Value = "Green Onion"
$ws.Range("B16").Value = "Scallion - Fresh"
$ws.Range("C16").Value = "'1"
$ws.Range("D16").Value = "'22.50"
$ws.Range("E16").Value = "'22.50"

# Row 17 - Sprouts - Alfalfa
$ws.Range("A17").Value = "'"
$ws.Range("B17").Value = "Sprouts - Alfalfa"
$ws.Range("C17").Value = "'4"
$ws.Range("D17").Value = "'16.85"
$ws.Range("E17").Value = "'67.40"

# Row 18 - Tomato - Fresh 5x6 Vine Ripe
$ws.Range("A18").Value = "'"
$ws.Range("B18").Value = "Tomato - Fresh 5x6 Vine Ripe"
$ws.Range("C18").Value = "'2"
$ws.Range("D18").Value = "'22.85"
$ws.Range("E18").Value = "'45.70"

# Row 19 - Tomato - Grape
$ws.Range("A19").Value = "'"
$ws.Range("B19").Value = "Tomato - Grape"
$ws.Range("C19").Value = "'4"
$ws.Range("D19").Value = "'14.85"
$ws.Range("E19").Value = "'59.40"
